$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns D values that look numeric must be forced to text so Excel
# keeps them as inline/shared strings rather than converting to numbers.
$numericLookingRows = 5,7,17,18,21,23,25,27,30,34,40,42,44,46,49,50,51
foreach ($r in $numericLookingRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = '29.100.41'
$ws.Cells.Item(2, 5).Value = '  +2.87%  '
$ws.Cells.Item(3, 4).Value = '1.582.46'
$ws.Cells.Item(3, 5).Value = '  +1.97%  '
$ws.Cells.Item(4, 5).Value = '  -0.67%  '
$ws.Cells.Item(5, 4).Value = '211.94'
$ws.Cells.Item(5, 5).Value = '  +1.34%  '
$ws.Cells.Item(6, 5).Value = '  +7.40%  '
$ws.Cells.Item(7, 4).Value = '0.994'
$ws.Cells.Item(7, 5).Value = '  -0.75%  '
$ws.Cells.Item(8, 5).Value = '  +9.02%  '
$ws.Cells.Item(9, 5).Value = '  +3.19%  '
$ws.Cells.Item(10, 5).Value = '  +1.84%  '
$ws.Cells.Item(11, 5).Value = '  +1.37%  '
$ws.Cells.Item(12, 4).Value = '1.809.52'
$ws.Cells.Item(12, 5).Value = '  +2.02%  '
$ws.Cells.Item(13, 4).Value = '1.569.02'
$ws.Cells.Item(13, 5).Value = '  +1.11%  '
$ws.Cells.Item(14, 4).Value = '29.139.42'
$ws.Cells.Item(14, 5).Value = '  +2.98%  '
$ws.Cells.Item(15, 5).Value = '  +2.78%  '
$ws.Cells.Item(16, 5).Value = '  +1.92%  '
$ws.Cells.Item(17, 4).Value = '62.38'
$ws.Cells.Item(17, 5).Value = '  +3.24%  '
$ws.Cells.Item(18, 4).Value = '238.55'
$ws.Cells.Item(18, 5).Value = '  +5.38%  '
$ws.Cells.Item(19, 5).Value = '  +1.96%  '
$ws.Cells.Item(20, 5).Value = '  +2.74%  '
$ws.Cells.Item(21, 4).Value = '0.993'
$ws.Cells.Item(21, 5).Value = '  -0.80%  '
$ws.Cells.Item(22, 5).Value = '  +2.31%  '
$ws.Cells.Item(23, 4).Value = '9.20'
$ws.Cells.Item(23, 5).Value = '  +4.51%  '
$ws.Cells.Item(24, 5).Value = '  +4.94%  '
$ws.Cells.Item(25, 4).Value = '152.86'
$ws.Cells.Item(25, 5).Value = '  +3.35%  '
$ws.Cells.Item(26, 5).Value = '  +5.02%  '
$ws.Cells.Item(27, 4).Value = '15.17'
$ws.Cells.Item(27, 5).Value = '  +2.76%  '
$ws.Cells.Item(28, 5).Value = '  +1.83%  '
$ws.Cells.Item(29, 5).Value = '  -0.72%  '
$ws.Cells.Item(30, 4).Value = '0.0464'
$ws.Cells.Item(30, 5).Value = '  -0.51%  '
$ws.Cells.Item(31, 5).Value = '  +0.04%  '
$ws.Cells.Item(32, 5).Value = '  +1.51%  '
$ws.Cells.Item(33, 4).Value = '1.423.76'
$ws.Cells.Item(33, 5).Value = '  +2.80%  '
$ws.Cells.Item(34, 4).Value = '3.05'
$ws.Cells.Item(34, 5).Value = '  -0.28%  '
$ws.Cells.Item(35, 5).Value = '  -1.28%  '
$ws.Cells.Item(36, 5).Value = '  +1.05%  '
$ws.Cells.Item(37, 5).Value = '  +7.37%  '
$ws.Cells.Item(38, 5).Value = '  -1.69%  '
$ws.Cells.Item(39, 5).Value = '  +1.29%  '
$ws.Cells.Item(40, 4).Value = '0.526'
$ws.Cells.Item(40, 5).Value = '  +2.74%  '
$ws.Cells.Item(41, 5).Value = '  +0.74%  '
$ws.Cells.Item(42, 4).Value = '0.994'
$ws.Cells.Item(42, 5).Value = '  -0.71%  '
$ws.Cells.Item(43, 5).Value = '  +1.53%  '
$ws.Cells.Item(44, 4).Value = '52.28'
$ws.Cells.Item(44, 5).Value = '  +24.60%  '
$ws.Cells.Item(45, 5).Value = '  -1.38%  '
$ws.Cells.Item(46, 4).Value = '64.75'
$ws.Cells.Item(46, 5).Value = '  +4.90%  '
$ws.Cells.Item(47, 5).Value = '  -1.75%  '
$ws.Cells.Item(48, 4).Value = '1.721.75'
$ws.Cells.Item(48, 5).Value = '  +2.04%  '
$ws.Cells.Item(49, 4).Value = '0.840'
$ws.Cells.Item(49, 5).Value = '  -7.30%  '
$ws.Cells.Item(50, 4).Value = '85.44'
$ws.Cells.Item(50, 5).Value = '  +0.01%  '
$ws.Cells.Item(51, 4).Value = '0.0512'
$ws.Cells.Item(51, 5).Value = '  +0.85%  '

# Restore the default cell style on the forced-text cells so no stray
# number-format style is left attached to them.
foreach ($r in $numericLookingRows) {
    $ws.Cells.Item($r, 4).Style = "Normal"
}
